# Edit script generated to apply the "take last response, not average" fix
# to Study2ValidationOnsetTimes.xlsx (Sheet1 rows 2-25, columns A:I, plus
# removal of the stray K2 SUM formula and a cosmetic selection change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2 held a leftover "=SUM(I:I)" debugging formula that is no longer needed.
$ws.Range("K2").ClearContents()

# Row 2
$ws.Range("C2").Value = 23.053829999975051
$ws.Range("F2").Value = 0.87382999997505095

# Row 3
$ws.Range("C3").Value = 35.32943333333025
$ws.Range("F3").Value = 0.89943333333025066
$ws.Range("I3").Value = 1
$ws.Range("A3:I3").Style = "Normal"

# Row 4
$ws.Range("C4").Value = 23.54635384616234
$ws.Range("F4").Value = 1.0963538461623408
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 92.857142857142861

# Row 5
$ws.Range("C5").Value = 16.110146153878915
$ws.Range("F5").Value = 2.1701461538789157

# Row 6
$ws.Range("C6").Value = 13.743663636365817
$ws.Range("F6").Value = 1.3236636363658167
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 78.571428571428569
$ws.Range("I6").Value = 0
$ws.Range("A6:I6").Style = "Bad"

# Row 7
$ws.Range("C7").Value = 29.981492307696872
$ws.Range("F7").Value = 0.26149230769687293
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 92.857142857142861

# Row 8
$ws.Range("C8").Value = 36.29034615386049
$ws.Range("F8").Value = 1.3803461538604935
$ws.Range("G8").Value = 13
$ws.Range("H8").Value = 92.857142857142861

# Row 9
$ws.Range("C9").Value = 43.283971428564392
$ws.Range("F9").Value = 1.3539714285643925

# Row 10
$ws.Range("C10").Value = 18.003321428594113
$ws.Range("F10").Value = 1.1833214285941125

# Row 11
$ws.Range("C11").Value = 32.869807142848934
$ws.Range("F11").Value = 1.1198071428489342

# Row 12
$ws.Range("C12").Value = 23.600914285722499
$ws.Range("F12").Value = 1.6009142857224994

# Row 13
$ws.Range("C13").Value = 26.943157142849355
$ws.Range("F13").Value = 0.34315714284935339
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 50

# Row 14
$ws.Range("C14").Value = 21.497992857157271
$ws.Range("F14").Value = 0.69799285715727066

# Row 15
$ws.Range("C15").Value = 11.176023076922972
$ws.Range("F15").Value = 1.586023076922972

# Row 17
$ws.Range("C17").Value = 13.596830769221093
$ws.Range("F17").Value = 0.81683076922109343

# Row 18
$ws.Range("C18").Value = 12.030841666658409
$ws.Range("F18").Value = 0.93084166665840939
$ws.Range("I18").Value = 1
$ws.Range("A18:I18").Style = "Normal"

# Row 19
$ws.Range("C19").Value = 29.166261538465264
$ws.Range("F19").Value = 1.2562615384652638
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 92.857142857142861

# Row 20
$ws.Range("C20").Value = 24.900800000001478
$ws.Range("F20").Value = 0.70080000000147891

# Row 21
$ws.Range("C21").Value = 9.260092857115394
$ws.Range("F21").Value = 1.260092857115394

# Row 22
$ws.Range("C22").Value = 24.442507142851639
$ws.Range("F22").Value = 1.7825071428516388

# Row 23
$ws.Range("C23").Value = 22.683685714289119
$ws.Range("F23").Value = 1.9636857142891202

# Row 24
$ws.Range("I24").Value = 1
$ws.Range("A24:I24").Style = "Normal"

# Row 25
$ws.Range("C25").Value = 36.657511111114154
$ws.Range("F25").Value = 1.177511111114157
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 64.285714285714292

# Restore the selection left behind by the author after the edit.
[void]$ws.Range("F27").Select()
